$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("books")

# C11/C12 used to hold placeholder junk text ("hhhhhhhhh" / "kjh;lkj;kl") in
# the ISBN column - replace them with a real (numeric) ISBN value, matching
# the ISBN already used for "wild" in row 9.
$ws.Cells.Item(11, 3).Value = 385486804
$ws.Cells.Item(12, 3).Value = 385486804

# Row 17 was a scratch/test row ("baking pi", rating 3.14, stray ISBN, rating
# "6") - clear it out, keeping only a real ISBN value in column C so the used
# range still reaches row 17.
$ws.Cells.Item(17, 1).Clear()
$ws.Cells.Item(17, 2).Clear()
$ws.Cells.Item(17, 4).Clear()
$ws.Cells.Item(17, 3).Value = 439136350

# Reset the view: scroll back to the top and move the selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C14").Select()
